$d = $word.ActiveDocument

# --- 1) Insert a new "Meta description" paragraph right after the title
#        heading paragraph. ---
$p1 = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)

# A single new <w:p> fragment passed to InsertXML merges into the following
# paragraph instead of creating a break, so a second (throwaway) <w:p/> is
# appended to force the paragraph split; it is deleted again right after.
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the unique adventure of Dragon Match, an online slot game with free spins, wilds, and cascading reels. Play for free today!</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($metaXml)

# Remove the spurious empty paragraph introduced above (it sits right after
# the new "Meta description" paragraph).
$emptyPara = $d.Paragraphs.Item(3)
$emptyPara.Range.Delete()

# --- 2) Remove the old duplicate "Play Dragon Match for Free..." paragraph
#        that used to sit right before the closing "Experience the unique
#        adventure..." paragraph. (It is plain/"Normal" styled, unlike the
#        similarly worded Heading 1 title at the very start of the doc, so
#        match on both text and style, searching from the end to be safe.) ---
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Text.TrimEnd() -eq "Play Dragon Match for Free - Exciting Online Slot Game") -and `
        ($p.Style.NameLocal -eq "Normal")) {
        $p.Range.Delete()
        break
    }
}

# --- 3) Replace the closing italic paragraph's text with the new feature
#        image prompt, scoping the Find/Replace to just that last paragraph
#        so the similarly-worded text inside the new "Meta description"
#        paragraph is left untouched. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute(
    "Experience the unique adventure of Dragon Match, an online slot game with free spins, wilds, and cascading reels. Play for free today!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create a feature image for Dragon Match: Please design a cartoon-style image featuring a happy Maya warrior with glasses. The warrior should be wearing traditional Maya clothing, including a colorful headband and jewelry. The background of the image should be inspired by Asian mythology with dragons and clouds. The warrior should also be holding a stack of gold coins and surrounded by cascading reels and game symbols. The overall style of the image should be fun and vibrant, capturing the adventurous spirit of the Dragon Match game.",
    2)
